$wb = $excel.ActiveWorkbook

# The three "Modify Transaction*" sheets each get a new row inserted before
# the final (submit) row, holding a ReceiptNumber label + a sheet-specific
# numeric value. The previously-last row shifts down by one.
$receiptValues = @{
    "Modify Transaction"  = 123
    "Modify Transaction1" = 543
    "Modify Transaction2" = 223
}

$selections = @{
    "Modify Transaction"  = "C11"
    "Modify Transaction1" = "C8"
    "Modify Transaction2" = "B11"
}

foreach ($name in @("Modify Transaction", "Modify Transaction1", "Modify Transaction2")) {
    $ws = $wb.Worksheets.Item($name)

    # Insert a new row 4 (pushes the existing row 4 down to row 5).
    $ws.Rows.Item(4).Insert()

    $ws.Cells.Item(4, 1).Value = "ReceiptNumber"
    $ws.Cells.Item(4, 2).Value = $receiptValues[$name]

    # Update the sheet's remembered selection without leaving it as the
    # active tab (restored below).
    [void]$ws.Range($selections[$name]).Select()
}

# Restore the originally active sheet/tab so tabSelected / activeTab are
# unchanged by our selection updates above.
[void]$wb.Worksheets.Item("Modify Transaction3").Activate()
